$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-09 (row 22)
$ws.Range("B22").Value = 6288
$ws.Range("D22").Value = 5837767
$ws.Range("E22").Value = 928.3980597964377
$ws.Range("F22").Value = 8.245825443277678
$ws.Range("H22").Value = 26.95323985072939
